$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column O: CONCAT(UPPER(D),"(""",D,"""),") for rows 1-122
$ws.Range("O1").Formula = '=CONCAT(UPPER(D1),"(""",D1,"""),")'
$ws.Range("O2:O65").Formula = '=CONCAT(UPPER(D2),"(""",D2,"""),")'
$ws.Range("O66:O122").Formula = '=CONCAT(UPPER(D66),"(""",D66,"""),")'

# Update view: scroll so column B is the left-most visible column,
# and select column O (mirrors the original selecting column J)
$ws.Range("O1").Select()
$excel.ActiveWindow.ScrollColumn = 2
